# Fetch Data Based on Col Name & Suite Runnable Check Scenario Updated
#
# This script:
#  1) Updates the selection on "TestData" to A1:B3
#  2) Updates the selection on "TestDataSet" to C3 (and makes it lose the
#     previous tabSelected flag since focus moves on)
#  3) Adds two new worksheets at the end of the workbook: "TestCase" (left
#     empty) and "TestSuite" (populated with suite-runnable data)
#  4) Formats the new "TestSuite" table the same way the existing tables in
#     the workbook are formatted (header row style + bordered data rows),
#     autofits column A, and finishes with B3 selected - making TestSuite
#     the active (last-visited) sheet/tab.

$wb = $excel.ActiveWorkbook

$testData    = $wb.Worksheets.Item(1)   # "TestData"
$testDataSet = $wb.Worksheets.Item(2)   # "TestDataSet"

# --- 1) Re-select a range on TestData -------------------------------------
$testData.Select()
$testData.Range("A1:B3").Select()

# --- 2) Re-select a cell on TestDataSet ------------------------------------
$testDataSet.Select()
$testDataSet.Range("C3").Select()

# --- 3) Add the "TestCase" worksheet (kept empty) --------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$testCase = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$testCase.Name = "TestCase"

# --- 4) Add the "TestSuite" worksheet and populate it -----------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$testSuite = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$testSuite.Name = "TestSuite"

$testSuite.Range("A1").Value = "Suite"
$testSuite.Range("B1").Value = "RunMode"
$testSuite.Range("A2").Value = "BankManagerSuite"
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("A3").Value = "CustomerSuite"
$testSuite.Range("B3").Value = "N"

# Match the formatting already used for tables elsewhere in the workbook:
# header row style, and the bordered "data row" style for the rest.
$testData.Range("A1:B1").Copy()
$testSuite.Range("A1:B1").PasteSpecial(-4122)

$testData.Range("A2:B2").Copy()
$testSuite.Range("A2:B3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Autofit column A like the other data sheets (bestFit columns).
$testSuite.Columns.Item(1).AutoFit()

# Leave TestSuite as the active sheet with B3 selected.
$testSuite.Select()
$testSuite.Range("B3").Select()
